# Adapt column header formatting to respective input file names (#7)
# - rename "<Header>_old" -> "<Header>_FV2410"
# - rename "<Header>_new" -> "<Header>_FV2504"
# - turn the data range into an Excel Table (ListObject)
# - freeze the header row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the used range extent (header row is row 1).
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count
$lastCol = $usedRange.Columns.Count

# Rename header cells: "_old" suffix -> "_FV2410", "_new" suffix -> "_FV2504".
for ($c = 1; $c -le $lastCol; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $header = $cell.Value2
    if ($header -ne $null) {
        if ($header.EndsWith("_old")) {
            $newHeader = $header.Substring(0, $header.Length - 4) + "_FV2410"
            $cell.Value = $newHeader
        } elseif ($header.EndsWith("_new")) {
            $newHeader = $header.Substring(0, $header.Length - 4) + "_FV2504"
            $cell.Value = $newHeader
        }
    }
}

# Turn the range into a proper Excel table ("Table1") covering the full
# used range, with the (now renamed) first row as the header row.
$tableRange = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item($lastRow, $lastCol))
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $tableRange, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

# Freeze the header row (split below row 1, top-left cell of the scrolling
# pane is A2).
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
